# Scheduled market-data refresh: update currentAveragePrice* / Leve*Price* / Leve*Profit*
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3100227
$ws.Range("I43").Value = 5852018.5
$ws.Range("J43").Value = 4461.5
$ws.Range("K43").Value = 5852018.5
$ws.Range("L43").Value = 4461.5
$ws.Range("M43").Value = -5851949.5
$ws.Range("N43").Value = -4599.5

$ws.Range("H86").Value = 1384.4054
$ws.Range("J86").Value = 1640.5
$ws.Range("L86").Value = 1640.5
$ws.Range("N86").Value = -3886.5

$ws.Range("H87").Value = 7223180
$ws.Range("I87").Value = 58000
$ws.Range("J87").Value = 8417377
$ws.Range("K87").Value = 58000
$ws.Range("L87").Value = 8417377
$ws.Range("M87").Value = -56752
$ws.Range("N87").Value = -8419873

$ws.Range("H89").Value = 1384.4054
$ws.Range("J89").Value = 1640.5
$ws.Range("L89").Value = 8202.5
$ws.Range("N89").Value = -19434.5

$ws.Range("H90").Value = 7223180
$ws.Range("I90").Value = 58000
$ws.Range("J90").Value = 8417377
$ws.Range("K90").Value = 174000
$ws.Range("L90").Value = 25252131
$ws.Range("M90").Value = -167760
$ws.Range("N90").Value = -25264611

$ws.Range("H99").Value = 354.77777
$ws.Range("I99").Value = 199
$ws.Range("J99").Value = 666.3333
$ws.Range("K99").Value = 597
$ws.Range("L99").Value = 1998.9999
$ws.Range("M99").Value = 901
$ws.Range("N99").Value = -4994.9999

$ws.Range("H112").Value = 627263.4399999999
$ws.Range("J112").Value = 771862.7
$ws.Range("L112").Value = 2315588.1
$ws.Range("N112").Value = -2317804.1

$ws.Range("H116").Value = 2647.5217
$ws.Range("I116").Value = 2268.2144
$ws.Range("J116").Value = 3237.5557
$ws.Range("K116").Value = 2268.2144
$ws.Range("L116").Value = 3237.5557
$ws.Range("M116").Value = 1173.7856
$ws.Range("N116").Value = -10121.5557

$ws.Range("H131").Value = 2906.2632
$ws.Range("I131").Value = 3417.0667
$ws.Range("J131").Value = 990.75
$ws.Range("K131").Value = 10251.2001
$ws.Range("L131").Value = 2972.25
$ws.Range("M131").Value = -5211.2001
$ws.Range("N131").Value = -13052.25

$ws.Range("H132").Value = 13515499
$ws.Range("I132").Value = 16951296
$ws.Range("J132").Value = 1363.4
$ws.Range("K132").Value = 50853888
$ws.Range("L132").Value = 4090.2
$ws.Range("M132").Value = -50851358
$ws.Range("N132").Value = -9150.200000000001

$ws.Range("H135").Value = 2217.4348
$ws.Range("I135").Value = 970.5714
$ws.Range("K135").Value = 8735.142600000001
$ws.Range("M135").Value = -6200.142600000001

$ws.Range("H138").Value = 160702.19
$ws.Range("I138").Value = 1313
$ws.Range("J138").Value = 194419.14
$ws.Range("K138").Value = 3939
$ws.Range("L138").Value = 583257.42
$ws.Range("M138").Value = 1201
$ws.Range("N138").Value = -593537.42

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11371.553
$ws.Range("I32").Value = 11371.553
$ws.Range("K32").Value = 11371.553
$ws.Range("M32").Value = -11084.553

$ws.Range("H122").Value = 4387.731
$ws.Range("I122").Value = 4191.6
$ws.Range("J122").Value = 4655.1816
$ws.Range("K122").Value = 12574.8
$ws.Range("L122").Value = 13965.5448
$ws.Range("M122").Value = -10124.8
$ws.Range("N122").Value = -18865.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 32709.766
$ws.Range("I82").Value = 6910.6665
$ws.Range("J82").Value = 61733.75
$ws.Range("K82").Value = 6910.6665
$ws.Range("L82").Value = 61733.75
$ws.Range("M82").Value = -6527.6665
$ws.Range("N82").Value = -62499.75

$ws.Range("H85").Value = 32709.766
$ws.Range("I85").Value = 6910.6665
$ws.Range("J85").Value = 61733.75
$ws.Range("K85").Value = 6910.6665
$ws.Range("L85").Value = 61733.75
$ws.Range("M85").Value = -5584.6665
$ws.Range("N85").Value = -64385.75

$ws.Range("H105").Value = 4990.5
$ws.Range("I105").Value = 3386.8
$ws.Range("K105").Value = 3386.8
$ws.Range("M105").Value = -1639.8

$ws.Range("H118").Value = 39998.668
$ws.Range("J118").Value = 39998.668
$ws.Range("L118").Value = 39998.668
$ws.Range("N118").Value = -43312.668

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H134").Value = 4755.852
$ws.Range("I134").Value = 1645.862
$ws.Range("J134").Value = 8363.440000000001
$ws.Range("K134").Value = 4937.586
$ws.Range("L134").Value = 25090.32
$ws.Range("M134").Value = -2402.586
$ws.Range("N134").Value = -30160.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2326.6667
$ws.Range("I31").Value = 2192.2
$ws.Range("J31").Value = 2999
$ws.Range("K31").Value = 2192.2
$ws.Range("L31").Value = 2999
$ws.Range("M31").Value = -1897.2
$ws.Range("N31").Value = -3589

$ws.Range("H34").Value = 2326.6667
$ws.Range("I34").Value = 2192.2
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 2192.2
$ws.Range("L34").Value = 2999
$ws.Range("M34").Value = -1990.2
$ws.Range("N34").Value = -3403

$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -9623
$ws.Range("N38").Value = ""

$ws.Range("H44").Value = 24750
$ws.Range("J44").Value = 24750
$ws.Range("L44").Value = 24750
$ws.Range("N44").Value = -25634

$ws.Range("H45").Value = 26000
$ws.Range("J45").Value = 26000
$ws.Range("L45").Value = 26000
$ws.Range("N45").Value = -27186

$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -9789
$ws.Range("N46").Value = ""

$ws.Range("H51").Value = 94998.75
$ws.Range("J51").Value = 94998.75
$ws.Range("L51").Value = 94998.75
$ws.Range("N51").Value = -96470.75

$ws.Range("H58").Value = 2979.3333
$ws.Range("J58").Value = 3099.6155
$ws.Range("L58").Value = 3099.6155
$ws.Range("N58").Value = -3505.6155

$ws.Range("H61").Value = 94998.75
$ws.Range("J61").Value = 94998.75
$ws.Range("L61").Value = 94998.75
$ws.Range("N61").Value = -95694.75

$ws.Range("H97").Value = 93299.336
$ws.Range("J97").Value = 93299.336
$ws.Range("L97").Value = 93299.336
$ws.Range("N97").Value = -95281.336

$ws.Range("H134").Value = 4902.5454
$ws.Range("I134").Value = 2941.4546
$ws.Range("J134").Value = 6863.636
$ws.Range("K134").Value = 8824.363799999999
$ws.Range("L134").Value = 20590.908
$ws.Range("M134").Value = -6289.363799999999
$ws.Range("N134").Value = -25660.908

$ws.Range("H136").Value = 2979.3333
$ws.Range("J136").Value = 3099.6155
$ws.Range("L136").Value = 9298.8465
$ws.Range("N136").Value = -14398.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2702.2632
$ws.Range("I132").Value = 1275.0625
$ws.Range("J132").Value = 3740.2273
$ws.Range("K132").Value = 11475.5625
$ws.Range("L132").Value = 33662.0457
$ws.Range("M132").Value = -8945.5625
$ws.Range("N132").Value = -38722.0457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 19989
$ws.Range("J35").Value = 19989
$ws.Range("L35").Value = 19989
$ws.Range("N35").Value = -20585

$ws.Range("H46").Value = 48263.57
$ws.Range("J46").Value = 49974.168
$ws.Range("L46").Value = 49974.168
$ws.Range("N46").Value = -50286.168

$ws.Range("H80").Value = 2470.2727
$ws.Range("I80").Value = 2237.25
$ws.Range("J80").Value = 2603.4285
$ws.Range("K80").Value = 2237.25
$ws.Range("L80").Value = 2603.4285
$ws.Range("M80").Value = -1239.25
$ws.Range("N80").Value = -4599.4285

$ws.Range("H83").Value = 2470.2727
$ws.Range("I83").Value = 2237.25
$ws.Range("J83").Value = 2603.4285
$ws.Range("K83").Value = 11186.25
$ws.Range("L83").Value = 13017.1425
$ws.Range("M83").Value = -6194.25
$ws.Range("N83").Value = -23001.1425

$ws.Range("H122").Value = 4244.815
$ws.Range("I122").Value = 3653.9443
$ws.Range("J122").Value = 5426.5557
$ws.Range("K122").Value = 10961.8329
$ws.Range("L122").Value = 16279.6671
$ws.Range("M122").Value = -8511.832900000001
$ws.Range("N122").Value = -21179.6671

$ws.Range("H132").Value = 13335968
$ws.Range("I132").Value = 13335968
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 40007904
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -40005374
$ws.Range("N132").Value = ""

$ws.Range("H137").Value = 135169.25
$ws.Range("J137").Value = 135169.25
$ws.Range("L137").Value = 135169.25
$ws.Range("N137").Value = -145369.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5762.7354
$ws.Range("I82").Value = 6840.731
$ws.Range("J82").Value = 2259.25
$ws.Range("K82").Value = 6840.731
$ws.Range("L82").Value = 2259.25
$ws.Range("M82").Value = -6479.731
$ws.Range("N82").Value = -2981.25

$ws.Range("H85").Value = 5762.7354
$ws.Range("I85").Value = 6840.731
$ws.Range("J85").Value = 2259.25
$ws.Range("K85").Value = 6840.731
$ws.Range("L85").Value = 2259.25
$ws.Range("M85").Value = -5592.731
$ws.Range("N85").Value = -4755.25

$ws.Range("H132").Value = 2875
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -13685

$ws.Range("H136").Value = 4559.968
$ws.Range("I136").Value = 4079.1538
$ws.Range("J136").Value = 7060.2
$ws.Range("K136").Value = 12237.4614
$ws.Range("L136").Value = 21180.6
$ws.Range("M136").Value = -9687.4614
$ws.Range("N136").Value = -26280.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1388.6666
$ws.Range("I100").Value = 1508.1666
$ws.Range("J100").Value = 1149.6666
$ws.Range("K100").Value = 3016.3332
$ws.Range("L100").Value = 2299.3332
$ws.Range("M100").Value = -2475.3332
$ws.Range("N100").Value = -3381.3332

$ws.Range("H113").Value = 2383587.8
$ws.Range("I113").Value = 3626588.8
$ws.Range("K113").Value = 10879766.4
$ws.Range("M113").Value = -10877596.4

$ws.Range("H136").Value = 7894.9375
$ws.Range("I136").Value = 9879.6
$ws.Range("J136").Value = 806.8570999999999
$ws.Range("K136").Value = 29638.8
$ws.Range("L136").Value = 2420.5713
$ws.Range("M136").Value = -27088.8
$ws.Range("N136").Value = -7520.5713
